# "nuevos experimentos no convexos"
# Refresh the generated-experiment numbers/expressions for a new run of the
# MitsosBarton2006Ex323 alpha-zero stationary-point generator.
#
# All of these cells store plain text in the original workbook (every value,
# numeric-looking or not, is a shared string) so purely numeric replacement
# values are written with the cell pre-formatted as Text ("@") to stop Excel
# from auto-converting them to numbers; expressions containing x/y are left
# alone since Excel already keeps those as text.

$wb = $excel.ActiveWorkbook

# NOTE: worksheet name lookups via Worksheets.Item(name) are case-insensitive
# and this workbook has both "Vector_bf" and "Vector_BF" sheets, so every
# sheet is addressed by its (1-based) tab position to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "1.9399999999999995 - x"
Set-TextValue $ws.Range("B2") "-2.9399999999999995"
Set-TextValue $ws.Range("D2") "0.32"
$ws.Range("A3").Value = "-1.9399999999999997 + x"
Set-TextValue $ws.Range("B3") "0.9399999999999997"
Set-TextValue $ws.Range("D3") "0.02"
$ws.Range("A4").Value = "35.63239999999999 + x - y - 9(x^2)"
Set-TextValue $ws.Range("B4") "-34.63239999999999"
Set-TextValue $ws.Range("D4") "0.44"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "-19.71359999999999 + (-0.5 + x)*(y^2)"
Set-TextValue $ws.Range("B2") "19.71359999999999"
Set-TextValue $ws.Range("D2") "0.36"
Set-TextValue $ws.Range("E2") "3.5"
Set-TextValue $ws.Range("F2") "0"
$ws.Range("A3").Value = "-3.6999999999999993 + y"
Set-TextValue $ws.Range("B3") "2.6999999999999993"
Set-TextValue $ws.Range("D3") "0.43"
Set-TextValue $ws.Range("E3") "3.2"
Set-TextValue $ws.Range("F3") "0"
$ws.Range("A4").Value = "-5.699999999999999 - y"
Set-TextValue $ws.Range("B4") "-4.699999999999999"
Set-TextValue $ws.Range("D4") "0.43"
Set-TextValue $ws.Range("E4") "2.1"
Set-TextValue $ws.Range("F4") "0"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "1.9399999999999997"
Set-TextValue $ws.Range("B2") "3.6999999999999993"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "-4.836159999999998"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-36.57019999999998"
Set-TextValue $ws.Range("A3") "-37.95599999999999"
